$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7433.28044715269
$ws.Range("C2").Value = 7101.97526660992
$ws.Range("E2").Value = 3254.96976263854
$ws.Range("F2").Value = 39.7060428853524

# Row 3
$ws.Range("B3").Value = 7161.2076590009
$ws.Range("C3").Value = 6789.488550861
$ws.Range("E3").Value = 3057.32713604807
$ws.Range("F3").Value = 259.450653621211

# Row 4
$ws.Range("B4").Value = 7019.85160074584
$ws.Range("C4").Value = 6615.2244321625
$ws.Range("E4").Value = 3185.77375792032
$ws.Range("F4").Value = 257.541591253451

# Row 5
$ws.Range("B5").Value = 7245.92370636213
$ws.Range("C5").Value = 6298.09073735816
$ws.Range("E5").Value = 3355.68194834902
$ws.Range("F5").Value = 251.407195237799

# Row 6
$ws.Range("B6").Value = 3023.93647427365
$ws.Range("C6").Value = 4924.11873122765
$ws.Range("E6").Value = 3559.04527240315
$ws.Range("F6").Value = 202.631833484617

# Row 7
$ws.Range("B7").Value = 3182.61929661698
$ws.Range("C7").Value = 4863.61544777032
$ws.Range("E7").Value = 3828.57377752232
$ws.Range("F7").Value = 211.341217720527

# Row 13
$ws.Range("E13").Value = 3794.11233978663
$ws.Range("F13").Value = 37.1877459985013

# Row 14
$ws.Range("E14").Value = 3794.11233978663
$ws.Range("F14").Value = 39.3463877788428

# Row 15
$ws.Range("E15").Value = 4682.6611643608
$ws.Range("F15").Value = 204.859291225045
